$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert the two new rows first (while row numbers below are still
#    at their original positions) so later A1-style addressing of the
#    untouched rows stays correct.
# ------------------------------------------------------------------
$ws.Rows.Item(11).Insert()   # new multiplication test (1,25*5)
$ws.Rows.Item(15).Insert()   # new division test (16,5/5,5)  -- this is
                              # "old" row 14 position, which after the
                              # first insert now sits one row lower

# ------------------------------------------------------------------
# 2. Fill in the content of the two new rows.
# ------------------------------------------------------------------
$ws.Range("A11").Value = "1. Запустить калькулятор`n2. Ввести число '1,25'`n3. Выбрать операцию '*'`n4. Ввести число '5'`n5. Нажать '='"
$ws.Range("B11").Value = "6.25"
$ws.Range("B11").WrapText = $true
$ws.Range("B11").NumberFormat = "@"
$ws.Rows.Item(11).RowHeight = 75

$ws.Range("A15").Value = "1. Запустить калькулятор`n2. Ввести число '16,5'`n3. Выбрать операцию '/'`n4. Ввести число '5,5'`n5. Нажать '='"
$ws.Range("B15").Value = 3
$ws.Rows.Item(15).RowHeight = 75

# ------------------------------------------------------------------
# 3. Append the brand-new last test case (power 0^299) as row 32.
# ------------------------------------------------------------------
$ws.Range("A32").Value = "1. Запустить калькулятор`n2. Ввести число '0'`n3. Нажать '^'`n4. Ввести число '299'`n5. Нажать '='"
$ws.Range("B32").Value = 0
$ws.Rows.Item(32).RowHeight = 75

# ------------------------------------------------------------------
# 4. B2 switches from numeric 66.45 to text "66,45".
# ------------------------------------------------------------------
$ws.Range("B2").Value = "66,45"
$ws.Range("B2").NumberFormat = "@"

# ------------------------------------------------------------------
# 5. Correct two values that display with slightly different rounding
#    in the supplemented sheet.
# ------------------------------------------------------------------
$ws.Range("B18").Value = 3.3332999999999999
$ws.Range("B24").Value = 3.3912

# ------------------------------------------------------------------
# 6. Apply the Text number format to the rest of column B (values are
#    already in place, so this only changes formatting/style, not the
#    stored value type).
# ------------------------------------------------------------------
$ws.Range("B3:B10").NumberFormat = "@"
$ws.Range("B12:B32").NumberFormat = "@"

# ------------------------------------------------------------------
# 7. Column B a touch wider.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 63.140625

# ------------------------------------------------------------------
# 8. Selection / view - scroll back to top, select B31.
# ------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("B31").Select()

# ------------------------------------------------------------------
# 9. Workbook-level changes: hidden _FilterDatabase defined name over
#    the full data range (Excel leaves this behind once a filter has
#    been toggled on the sheet).
# ------------------------------------------------------------------
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Лист1!`$A`$1:`$B`$32")
$n.Visible = $false
